# Auto-generated Excel COM-interop script to apply data refresh changes
# to the leve-profit tables across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 188  # H2
$ws.Cells.Item(2, 9).Value = 188  # I2
$ws.Cells.Item(2, 11).Value = 188  # K2
$ws.Cells.Item(2, 13).Value = -75  # M2
$ws.Cells.Item(4, 8).Value = 183  # H4
$ws.Cells.Item(4, 9).Value = 183  # I4
$ws.Cells.Item(4, 11).Value = 183  # K4
$ws.Cells.Item(4, 13).Value = -69  # M4
$ws.Cells.Item(11, 8).Value = 201.5  # H11
$ws.Cells.Item(11, 9).Value = 201.5  # I11
$ws.Cells.Item(11, 11).Value = 201.5  # K11
$ws.Cells.Item(11, 13).Value = -61.5  # M11
$ws.Cells.Item(17, 8).Value = 2903324.5  # H17
$ws.Cells.Item(17, 10).Value = 3003221.8  # J17
$ws.Cells.Item(17, 12).Value = 9009665.399999999  # L17
$ws.Cells.Item(17, 14).Value = -9010001.399999999  # N17
$ws.Cells.Item(19, 8).Value = 886.6111  # H19
$ws.Cells.Item(19, 9).Value = 952.1667  # I19
$ws.Cells.Item(19, 10).Value = 755.5  # J19
$ws.Cells.Item(19, 11).Value = 952.1667  # K19
$ws.Cells.Item(19, 12).Value = 755.5  # L19
$ws.Cells.Item(19, 13).Value = -777.1667  # M19
$ws.Cells.Item(19, 14).Value = -1105.5  # N19
$ws.Cells.Item(28, 8).Value = 529.8421  # H28
$ws.Cells.Item(28, 9).Value = 337.77777  # I28
$ws.Cells.Item(28, 11).Value = 337.77777  # K28
$ws.Cells.Item(28, 13).Value = 147.22223  # M28
$ws.Cells.Item(32, 8).Value = 26455.385  # H32
$ws.Cells.Item(32, 9).Value = 47038.5  # I32
$ws.Cells.Item(32, 10).Value = 8812.714  # J32
$ws.Cells.Item(32, 11).Value = 47038.5  # K32
$ws.Cells.Item(32, 12).Value = 8812.714  # L32
$ws.Cells.Item(32, 13).Value = -46712.5  # M32
$ws.Cells.Item(32, 14).Value = -9464.714  # N32
$ws.Cells.Item(33, 8).Value = 3004569.5  # H33
$ws.Cells.Item(33, 9).Value = 6757231.5  # I33
$ws.Cells.Item(33, 10).Value = 2440  # J33
$ws.Cells.Item(33, 11).Value = 6757231.5  # K33
$ws.Cells.Item(33, 12).Value = 2440  # L33
$ws.Cells.Item(33, 13).Value = -6757002.5  # M33
$ws.Cells.Item(33, 14).Value = -2898  # N33
$ws.Cells.Item(40, 8).Value = 1606.4166  # H40
$ws.Cells.Item(40, 9).Value = 1625.1818  # I40
$ws.Cells.Item(40, 10).Value = 1400  # J40
$ws.Cells.Item(40, 11).Value = 1625.1818  # K40
$ws.Cells.Item(40, 12).Value = 1400  # L40
$ws.Cells.Item(40, 13).Value = -1450.1818  # M40
$ws.Cells.Item(40, 14).Value = -1750  # N40
$ws.Cells.Item(41, 8).Value = 803.6  # H41
$ws.Cells.Item(41, 10).Value = 931.1667  # J41
$ws.Cells.Item(41, 12).Value = 931.1667  # L41
$ws.Cells.Item(41, 14).Value = -1811.1667  # N41
$ws.Cells.Item(43, 8).Value = 2988.7  # H43
$ws.Cells.Item(43, 10).Value = 2599.5  # J43
$ws.Cells.Item(43, 12).Value = 2599.5  # L43
$ws.Cells.Item(43, 14).Value = -2737.5  # N43
$ws.Cells.Item(51, 8).Value = 6733.0713  # H51
$ws.Cells.Item(51, 9).Value = 6040.6665  # I51
$ws.Cells.Item(51, 10).Value = 7979.4  # J51
$ws.Cells.Item(51, 11).Value = 6040.6665  # K51
$ws.Cells.Item(51, 12).Value = 7979.4  # L51
$ws.Cells.Item(51, 13).Value = -5556.6665  # M51
$ws.Cells.Item(51, 14).Value = -8947.4  # N51
$ws.Cells.Item(64, 8).Value = 4533.3335  # H64
$ws.Cells.Item(64, 10).Value = 5000  # J64
$ws.Cells.Item(64, 12).Value = 5000  # L64
$ws.Cells.Item(64, 14).Value = -5496  # N64
$ws.Cells.Item(67, 8).Value = 4533.3335  # H67
$ws.Cells.Item(67, 10).Value = 5000  # J67
$ws.Cells.Item(67, 12).Value = 5000  # L67
$ws.Cells.Item(67, 14).Value = -6716  # N67
$ws.Cells.Item(74, 8).Value = 7850  # H74
$ws.Cells.Item(74, 9).Value = 7850  # I74
$ws.Cells.Item(74, 11).Value = 7850  # K74
$ws.Cells.Item(74, 13).Value = -6914  # M74
$ws.Cells.Item(77, 8).Value = 7850  # H77
$ws.Cells.Item(77, 9).Value = 7850  # I77
$ws.Cells.Item(77, 11).Value = 39250  # K77
$ws.Cells.Item(77, 13).Value = -34570  # M77
$ws.Cells.Item(86, 8).Value = 8726  # H86
$ws.Cells.Item(86, 9).Value = 5000  # I86
$ws.Cells.Item(86, 10).Value = 9968  # J86
$ws.Cells.Item(86, 11).Value = 5000  # K86
$ws.Cells.Item(86, 12).Value = 9968  # L86
$ws.Cells.Item(86, 13).Value = -3877  # M86
$ws.Cells.Item(86, 14).Value = -12214  # N86
$ws.Cells.Item(89, 8).Value = 8726  # H89
$ws.Cells.Item(89, 9).Value = 5000  # I89
$ws.Cells.Item(89, 10).Value = 9968  # J89
$ws.Cells.Item(89, 11).Value = 25000  # K89
$ws.Cells.Item(89, 12).Value = 49840  # L89
$ws.Cells.Item(89, 13).Value = -19384  # M89
$ws.Cells.Item(89, 14).Value = -61072  # N89
$ws.Cells.Item(92, 8).Value = 986.6  # H92
$ws.Cells.Item(92, 9).Value = 1027.625  # I92
$ws.Cells.Item(92, 11).Value = 1027.625  # K92
$ws.Cells.Item(92, 13).Value = 220.375  # M92
$ws.Cells.Item(96, 8).Value = 462.22223  # H96
$ws.Cells.Item(96, 9).Value = 352.375  # I96
$ws.Cells.Item(96, 11).Value = 1057.125  # K96
$ws.Cells.Item(96, 13).Value = 315.875  # M96
$ws.Cells.Item(100, 8).Value = 66068.19  # H100
$ws.Cells.Item(100, 9).Value = 80339.84  # I100
$ws.Cells.Item(100, 10).Value = 4224.3335  # J100
$ws.Cells.Item(100, 11).Value = 80339.84  # K100
$ws.Cells.Item(100, 12).Value = 4224.3335  # L100
$ws.Cells.Item(100, 13).Value = -79798.84  # M100
$ws.Cells.Item(100, 14).Value = -5306.3335  # N100
$ws.Cells.Item(106, 8).Value = 24885.32  # H106
$ws.Cells.Item(106, 10).Value = 28855.643  # J106
$ws.Cells.Item(106, 12).Value = 28855.643  # L106
$ws.Cells.Item(106, 14).Value = -30117.643  # N106
$ws.Cells.Item(107, 8).Value = 1153.2609  # H107
$ws.Cells.Item(107, 9).Value = 909.3  # I107
$ws.Cells.Item(107, 10).Value = 2779.6667  # J107
$ws.Cells.Item(107, 11).Value = 909.3  # K107
$ws.Cells.Item(107, 12).Value = 2779.6667  # L107
$ws.Cells.Item(107, 13).Value = 1010.7  # M107
$ws.Cells.Item(107, 14).Value = -6619.6667  # N107
$ws.Cells.Item(111, 8).Value = 1514  # H111
$ws.Cells.Item(111, 9).Value = 1514  # I111
$ws.Cells.Item(111, 11).Value = 4542  # K111
$ws.Cells.Item(111, 13).Value = -1475  # M111
$ws.Cells.Item(112, 8).Value = 1480.8214  # H112
$ws.Cells.Item(112, 10).Value = 1566.3043  # J112
$ws.Cells.Item(112, 12).Value = 4698.9129  # L112
$ws.Cells.Item(112, 14).Value = -6914.9129  # N112
$ws.Cells.Item(115, 8).Value = 1750  # H115
$ws.Cells.Item(115, 9).Value = 1000  # I115
$ws.Cells.Item(115, 11).Value = 3000  # K115
$ws.Cells.Item(115, 13).Value = -1433  # M115
$ws.Cells.Item(116, 8).Value = 12014.368  # H116
$ws.Cells.Item(116, 9).Value = 9774  # I116
$ws.Cells.Item(116, 11).Value = 9774  # K116
$ws.Cells.Item(116, 13).Value = -6332  # M116
$ws.Cells.Item(118, 8).Value = 900  # H118
$ws.Cells.Item(118, 9).Value = 900  # I118
$ws.Cells.Item(118, 11).Value = 2700  # K118
$ws.Cells.Item(118, 13).Value = -1043  # M118
$ws.Cells.Item(125, 8).Value = 1515  # H125
$ws.Cells.Item(125, 10).Value = 2499.5  # J125
$ws.Cells.Item(125, 12).Value = 22495.5  # L125
$ws.Cells.Item(125, 14).Value = -27415.5  # N125
$ws.Cells.Item(137, 8).Value = 12156.538  # H137
$ws.Cells.Item(137, 9).Value = 4831.636  # I137
$ws.Cells.Item(137, 10).Value = 21635.824  # J137
$ws.Cells.Item(137, 11).Value = 14494.908  # K137
$ws.Cells.Item(137, 12).Value = 64907.472  # L137
$ws.Cells.Item(137, 13).Value = -11944.908  # M137
$ws.Cells.Item(137, 14).Value = -70007.47200000001  # N137
$ws.Cells.Item(138, 8).Value = 1999.8695  # H138
$ws.Cells.Item(138, 10).Value = 0  # J138
$ws.Cells.Item(138, 12).Value = 0  # L138
$ws.Cells.Item(138, 14).ClearContents()  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2452.4714  # H32
$ws.Cells.Item(32, 9).Value = 1333.6936  # I32
$ws.Cells.Item(32, 11).Value = 1333.6936  # K32
$ws.Cells.Item(32, 13).Value = -1046.6936  # M32
$ws.Cells.Item(45, 8).Value = 14765.9  # H45
$ws.Cells.Item(45, 9).Value = 20263.5  # I45
$ws.Cells.Item(45, 10).Value = 6519.5  # J45
$ws.Cells.Item(45, 11).Value = 20263.5  # K45
$ws.Cells.Item(45, 12).Value = 6519.5  # L45
$ws.Cells.Item(45, 13).Value = -19886.5  # M45
$ws.Cells.Item(45, 14).Value = -7273.5  # N45
$ws.Cells.Item(46, 8).Value = 15294.667  # H46
$ws.Cells.Item(46, 9).Value = 13624.75  # I46
$ws.Cells.Item(46, 11).Value = 13624.75  # K46
$ws.Cells.Item(46, 13).Value = -13305.75  # M46
$ws.Cells.Item(61, 8).Value = 4773.207  # H61
$ws.Cells.Item(61, 9).Value = 3049.6191  # I61
$ws.Cells.Item(61, 10).Value = 9297.625  # J61
$ws.Cells.Item(61, 11).Value = 3049.6191  # K61
$ws.Cells.Item(61, 12).Value = 9297.625  # L61
$ws.Cells.Item(61, 13).Value = -2837.6191  # M61
$ws.Cells.Item(61, 14).Value = -9721.625  # N61
$ws.Cells.Item(63, 8).Value = 3837.3572  # H63
$ws.Cells.Item(63, 9).Value = 2934.75  # I63
$ws.Cells.Item(63, 11).Value = 2934.75  # K63
$ws.Cells.Item(63, 13).Value = -2248.75  # M63
$ws.Cells.Item(66, 8).Value = 3837.3572  # H66
$ws.Cells.Item(66, 9).Value = 2934.75  # I66
$ws.Cells.Item(66, 11).Value = 14673.75  # K66
$ws.Cells.Item(66, 13).Value = -11241.75  # M66
$ws.Cells.Item(74, 8).Value = 2496.6667  # H74
$ws.Cells.Item(74, 9).Value = 1846.4445  # I74
$ws.Cells.Item(74, 10).Value = 3667.0667  # J74
$ws.Cells.Item(74, 11).Value = 1846.4445  # K74
$ws.Cells.Item(74, 12).Value = 3667.0667  # L74
$ws.Cells.Item(74, 13).Value = -972.4445000000001  # M74
$ws.Cells.Item(74, 14).Value = -5415.066699999999  # N74
$ws.Cells.Item(77, 8).Value = 2496.6667  # H77
$ws.Cells.Item(77, 9).Value = 1846.4445  # I77
$ws.Cells.Item(77, 10).Value = 3667.0667  # J77
$ws.Cells.Item(77, 11).Value = 9232.2225  # K77
$ws.Cells.Item(77, 12).Value = 18335.3335  # L77
$ws.Cells.Item(77, 13).Value = -4864.2225  # M77
$ws.Cells.Item(77, 14).Value = -27071.3335  # N77
$ws.Cells.Item(136, 8).Value = 4773.207  # H136
$ws.Cells.Item(136, 9).Value = 3049.6191  # I136
$ws.Cells.Item(136, 10).Value = 9297.625  # J136
$ws.Cells.Item(136, 11).Value = 9148.8573  # K136
$ws.Cells.Item(136, 12).Value = 27892.875  # L136
$ws.Cells.Item(136, 13).Value = -6598.8573  # M136
$ws.Cells.Item(136, 14).Value = -32992.875  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 182.71428  # H80
$ws.Cells.Item(80, 10).Value = 58.666668  # J80
$ws.Cells.Item(80, 12).Value = 58.666668  # L80
$ws.Cells.Item(80, 14).Value = -2054.666668  # N80
$ws.Cells.Item(83, 8).Value = 182.71428  # H83
$ws.Cells.Item(83, 10).Value = 58.666668  # J83
$ws.Cells.Item(83, 12).Value = 293.33334  # L83
$ws.Cells.Item(83, 14).Value = -10277.33334  # N83
$ws.Cells.Item(86, 8).Value = 314477.3  # H86
$ws.Cells.Item(86, 10).Value = 2452.1  # J86
$ws.Cells.Item(86, 12).Value = 2452.1  # L86
$ws.Cells.Item(86, 14).Value = -4698.1  # N86
$ws.Cells.Item(89, 8).Value = 314477.3  # H89
$ws.Cells.Item(89, 10).Value = 2452.1  # J89
$ws.Cells.Item(89, 12).Value = 12260.5  # L89
$ws.Cells.Item(89, 14).Value = -23492.5  # N89
$ws.Cells.Item(99, 8).Value = 4400.1904  # H99
$ws.Cells.Item(99, 9).Value = 4444.3125  # I99
$ws.Cells.Item(99, 10).Value = 4259  # J99
$ws.Cells.Item(99, 11).Value = 4444.3125  # K99
$ws.Cells.Item(99, 12).Value = 4259  # L99
$ws.Cells.Item(99, 13).Value = -2946.3125  # M99
$ws.Cells.Item(99, 14).Value = -7255  # N99
$ws.Cells.Item(105, 8).Value = 3488.7368  # H105
$ws.Cells.Item(105, 9).Value = 3812.5  # I105
$ws.Cells.Item(105, 11).Value = 3812.5  # K105
$ws.Cells.Item(105, 13).Value = -2065.5  # M105
$ws.Cells.Item(107, 8).Value = 1087.0667  # H107
$ws.Cells.Item(107, 9).Value = 950.1923  # I107
$ws.Cells.Item(107, 11).Value = 950.1923  # K107
$ws.Cells.Item(107, 13).Value = 969.8077  # M107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 426.8  # H22
$ws.Cells.Item(22, 9).Value = 462.16666  # I22
$ws.Cells.Item(22, 10).Value = 373.75  # J22
$ws.Cells.Item(22, 11).Value = 462.16666  # K22
$ws.Cells.Item(22, 12).Value = 373.75  # L22
$ws.Cells.Item(22, 13).Value = -112.16666  # M22
$ws.Cells.Item(22, 14).Value = -1073.75  # N22
$ws.Cells.Item(31, 8).Value = 3780.2083  # H31
$ws.Cells.Item(31, 9).Value = 1677.3334  # I31
$ws.Cells.Item(31, 10).Value = 5883.0835  # J31
$ws.Cells.Item(31, 11).Value = 1677.3334  # K31
$ws.Cells.Item(31, 12).Value = 5883.0835  # L31
$ws.Cells.Item(31, 13).Value = -1382.3334  # M31
$ws.Cells.Item(31, 14).Value = -6473.0835  # N31
$ws.Cells.Item(34, 8).Value = 3780.2083  # H34
$ws.Cells.Item(34, 9).Value = 1677.3334  # I34
$ws.Cells.Item(34, 10).Value = 5883.0835  # J34
$ws.Cells.Item(34, 11).Value = 1677.3334  # K34
$ws.Cells.Item(34, 12).Value = 5883.0835  # L34
$ws.Cells.Item(34, 13).Value = -1475.3334  # M34
$ws.Cells.Item(34, 14).Value = -6287.0835  # N34
$ws.Cells.Item(62, 8).Value = 2817.3333  # H62
$ws.Cells.Item(62, 9).Value = 2899.6  # I62
$ws.Cells.Item(62, 10).Value = 2406  # J62
$ws.Cells.Item(62, 11).Value = 2899.6  # K62
$ws.Cells.Item(62, 12).Value = 2406  # L62
$ws.Cells.Item(62, 13).Value = -2275.6  # M62
$ws.Cells.Item(62, 14).Value = -3654  # N62
$ws.Cells.Item(65, 8).Value = 2817.3333  # H65
$ws.Cells.Item(65, 9).Value = 2899.6  # I65
$ws.Cells.Item(65, 10).Value = 2406  # J65
$ws.Cells.Item(65, 11).Value = 14498  # K65
$ws.Cells.Item(65, 12).Value = 12030  # L65
$ws.Cells.Item(65, 13).Value = -11378  # M65
$ws.Cells.Item(65, 14).Value = -18270  # N65
$ws.Cells.Item(93, 8).Value = 23724.25  # H93
$ws.Cells.Item(93, 9).Value = 21666.334  # I93
$ws.Cells.Item(93, 10).Value = 29898  # J93
$ws.Cells.Item(93, 11).Value = 21666.334  # K93
$ws.Cells.Item(93, 12).Value = 29898  # L93
$ws.Cells.Item(93, 13).Value = -19794.334  # M93
$ws.Cells.Item(93, 14).Value = -33642  # N93
$ws.Cells.Item(105, 8).Value = 1458.6666  # H105
$ws.Cells.Item(105, 9).Value = 1574.625  # I105
$ws.Cells.Item(105, 10).Value = 531  # J105
$ws.Cells.Item(105, 11).Value = 1574.625  # K105
$ws.Cells.Item(105, 12).Value = 531  # L105
$ws.Cells.Item(105, 13).Value = 172.375  # M105
$ws.Cells.Item(105, 14).Value = -4025  # N105
$ws.Cells.Item(107, 8).Value = 693.6923  # H107
$ws.Cells.Item(107, 9).Value = 580.64  # I107
$ws.Cells.Item(107, 10).Value = 895.5714  # J107
$ws.Cells.Item(107, 11).Value = 580.64  # K107
$ws.Cells.Item(107, 12).Value = 895.5714  # L107
$ws.Cells.Item(107, 13).Value = 1339.36  # M107
$ws.Cells.Item(107, 14).Value = -4735.5714  # N107
$ws.Cells.Item(122, 8).Value = 3007.4167  # H122
$ws.Cells.Item(122, 9).Value = 2772  # I122
$ws.Cells.Item(122, 10).Value = 3478.25  # J122
$ws.Cells.Item(122, 11).Value = 8316  # K122
$ws.Cells.Item(122, 12).Value = 10434.75  # L122
$ws.Cells.Item(122, 13).Value = -5866  # M122
$ws.Cells.Item(122, 14).Value = -15334.75  # N122

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 2943.2  # H14
$ws.Cells.Item(14, 9).Value = 2943.2  # I14
$ws.Cells.Item(14, 11).Value = 8829.599999999999  # K14
$ws.Cells.Item(14, 13).Value = -8656.599999999999  # M14
$ws.Cells.Item(46, 8).Value = 11500200  # H46
$ws.Cells.Item(46, 9).Value = 13800120  # I46
$ws.Cells.Item(46, 10).Value = 600  # J46
$ws.Cells.Item(46, 11).Value = 41400360  # K46
$ws.Cells.Item(46, 12).Value = 1800  # L46
$ws.Cells.Item(46, 13).Value = -41400269  # M46
$ws.Cells.Item(46, 14).Value = -1982  # N46
$ws.Cells.Item(50, 8).Value = 2156.8572  # H50
$ws.Cells.Item(50, 9).Value = 867.6667  # I50
$ws.Cells.Item(50, 11).Value = 2603.0001  # K50
$ws.Cells.Item(50, 13).Value = -2122.0001  # M50
$ws.Cells.Item(53, 8).Value = 2156.8572  # H53
$ws.Cells.Item(53, 9).Value = 867.6667  # I53
$ws.Cells.Item(53, 11).Value = 2603.0001  # K53
$ws.Cells.Item(53, 13).Value = -2122.0001  # M53
$ws.Cells.Item(69, 8).Value = 4332.4546  # H69
$ws.Cells.Item(69, 9).Value = 2536.875  # I69
$ws.Cells.Item(69, 10).Value = 9120.666999999999  # J69
$ws.Cells.Item(69, 11).Value = 7610.625  # K69
$ws.Cells.Item(69, 12).Value = 27362.001  # L69
$ws.Cells.Item(69, 13).Value = -6799.625  # M69
$ws.Cells.Item(69, 14).Value = -28984.001  # N69
$ws.Cells.Item(72, 8).Value = 4332.4546  # H72
$ws.Cells.Item(72, 9).Value = 2536.875  # I72
$ws.Cells.Item(72, 10).Value = 9120.666999999999  # J72
$ws.Cells.Item(72, 11).Value = 22831.875  # K72
$ws.Cells.Item(72, 12).Value = 82086.003  # L72
$ws.Cells.Item(72, 13).Value = -18775.875  # M72
$ws.Cells.Item(72, 14).Value = -90198.003  # N72
$ws.Cells.Item(114, 8).Value = 1914.1666  # H114
$ws.Cells.Item(114, 9).Value = 538.6667  # I114
$ws.Cells.Item(114, 10).Value = 2372.6667  # J114
$ws.Cells.Item(114, 11).Value = 1616.0001  # K114
$ws.Cells.Item(114, 12).Value = 7118.000100000001  # L114
$ws.Cells.Item(114, 13).Value = 1637.9999  # M114
$ws.Cells.Item(114, 14).Value = -13626.0001  # N114
$ws.Cells.Item(117, 8).Value = 909.5  # H117
$ws.Cells.Item(117, 10).Value = 544.75  # J117
$ws.Cells.Item(117, 12).Value = 1634.25  # L117
$ws.Cells.Item(117, 14).Value = -8518.25  # N117
$ws.Cells.Item(136, 8).Value = 2164.3333  # H136
$ws.Cells.Item(136, 10).Value = 1996.5  # J136
$ws.Cells.Item(136, 12).Value = 5989.5  # L136
$ws.Cells.Item(136, 14).Value = -16189.5  # N136

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 807.9545000000001  # H97
$ws.Cells.Item(97, 9).Value = 849.61536  # I97
$ws.Cells.Item(97, 10).Value = 747.7778  # J97
$ws.Cells.Item(97, 11).Value = 849.61536  # K97
$ws.Cells.Item(97, 12).Value = 747.7778  # L97
$ws.Cells.Item(97, 13).Value = -353.61536  # M97
$ws.Cells.Item(97, 14).Value = -1739.7778  # N97
$ws.Cells.Item(122, 8).Value = 951.5  # H122
$ws.Cells.Item(122, 9).Value = 933.8333  # I122
$ws.Cells.Item(122, 10).Value = 978  # J122
$ws.Cells.Item(122, 11).Value = 2801.4999  # K122
$ws.Cells.Item(122, 12).Value = 2934  # L122
$ws.Cells.Item(122, 13).Value = -351.4998999999998  # M122
$ws.Cells.Item(122, 14).Value = -7834  # N122
$ws.Cells.Item(129, 8).Value = 45000  # H129
$ws.Cells.Item(129, 10).Value = 45000  # J129
$ws.Cells.Item(129, 12).Value = 45000  # L129
$ws.Cells.Item(129, 14).Value = -55000  # N129

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4251  # H7
$ws.Cells.Item(7, 9).Value = 3797.1428  # I7
$ws.Cells.Item(7, 11).Value = 3797.1428  # K7
$ws.Cells.Item(7, 13).Value = -3685.1428  # M7
$ws.Cells.Item(22, 8).Value = 1791.9333  # H22
$ws.Cells.Item(22, 9).Value = 3182.3333  # I22
$ws.Cells.Item(22, 10).Value = 1444.3334  # J22
$ws.Cells.Item(22, 11).Value = 3182.3333  # K22
$ws.Cells.Item(22, 12).Value = 1444.3334  # L22
$ws.Cells.Item(22, 13).Value = -2887.3333  # M22
$ws.Cells.Item(22, 14).Value = -2034.3334  # N22
$ws.Cells.Item(27, 8).Value = 1791.9333  # H27
$ws.Cells.Item(27, 9).Value = 3182.3333  # I27
$ws.Cells.Item(27, 10).Value = 1444.3334  # J27
$ws.Cells.Item(27, 11).Value = 3182.3333  # K27
$ws.Cells.Item(27, 12).Value = 1444.3334  # L27
$ws.Cells.Item(27, 13).Value = -3075.3333  # M27
$ws.Cells.Item(27, 14).Value = -1658.3334  # N27
$ws.Cells.Item(40, 8).Value = 1900  # H40
$ws.Cells.Item(40, 9).Value = 1900  # I40
$ws.Cells.Item(40, 11).Value = 1900  # K40
$ws.Cells.Item(40, 13).Value = -1764  # M40
$ws.Cells.Item(46, 8).Value = 881.1539  # H46
$ws.Cells.Item(46, 9).Value = 1041.4  # I46
$ws.Cells.Item(46, 10).Value = 781  # J46
$ws.Cells.Item(46, 11).Value = 1041.4  # K46
$ws.Cells.Item(46, 12).Value = 781  # L46
$ws.Cells.Item(46, 13).Value = -853.4000000000001  # M46
$ws.Cells.Item(46, 14).Value = -1157  # N46
$ws.Cells.Item(55, 8).Value = 2090.923  # H55
$ws.Cells.Item(55, 9).Value = 2340.4285  # I55
$ws.Cells.Item(55, 10).Value = 1799.8334  # J55
$ws.Cells.Item(55, 11).Value = 2340.4285  # K55
$ws.Cells.Item(55, 12).Value = 1799.8334  # L55
$ws.Cells.Item(55, 13).Value = -2167.4285  # M55
$ws.Cells.Item(55, 14).Value = -2145.8334  # N55
$ws.Cells.Item(61, 8).Value = 2236.889  # H61
$ws.Cells.Item(61, 9).Value = 2110.6667  # I61
$ws.Cells.Item(61, 11).Value = 2110.6667  # K61
$ws.Cells.Item(61, 13).Value = -1908.6667  # M61
$ws.Cells.Item(68, 8).Value = 2560.3462  # H68
$ws.Cells.Item(68, 9).Value = 2586.0435  # I68
$ws.Cells.Item(68, 10).Value = 2363.3333  # J68
$ws.Cells.Item(68, 11).Value = 2586.0435  # K68
$ws.Cells.Item(68, 12).Value = 2363.3333  # L68
$ws.Cells.Item(68, 13).Value = -1837.0435  # M68
$ws.Cells.Item(68, 14).Value = -3861.3333  # N68
$ws.Cells.Item(71, 8).Value = 2560.3462  # H71
$ws.Cells.Item(71, 9).Value = 2586.0435  # I71
$ws.Cells.Item(71, 10).Value = 2363.3333  # J71
$ws.Cells.Item(71, 11).Value = 12930.2175  # K71
$ws.Cells.Item(71, 12).Value = 11816.6665  # L71
$ws.Cells.Item(71, 13).Value = -9186.217500000001  # M71
$ws.Cells.Item(71, 14).Value = -19304.6665  # N71
$ws.Cells.Item(80, 8).Value = 44276.4  # H80
$ws.Cells.Item(80, 9).Value = 0  # I80
$ws.Cells.Item(80, 10).Value = 44276.4  # J80
$ws.Cells.Item(80, 11).Value = 0  # K80
$ws.Cells.Item(80, 12).Value = 44276.4  # L80
$ws.Cells.Item(80, 13).ClearContents()  # M80
$ws.Cells.Item(80, 14).Value = -46522.4  # N80
$ws.Cells.Item(83, 8).Value = 44276.4  # H83
$ws.Cells.Item(83, 9).Value = 0  # I83
$ws.Cells.Item(83, 10).Value = 44276.4  # J83
$ws.Cells.Item(83, 11).Value = 0  # K83
$ws.Cells.Item(83, 12).Value = 132829.2  # L83
$ws.Cells.Item(83, 13).ClearContents()  # M83
$ws.Cells.Item(83, 14).Value = -144061.2  # N83
$ws.Cells.Item(93, 8).Value = 3722.05  # H93
$ws.Cells.Item(93, 9).Value = 3760.0527  # I93
$ws.Cells.Item(93, 11).Value = 3760.0527  # K93
$ws.Cells.Item(93, 13).Value = -2512.0527  # M93
$ws.Cells.Item(100, 8).Value = 374283.84  # H100
$ws.Cells.Item(100, 9).Value = 4104.913  # I100
$ws.Cells.Item(100, 11).Value = 4104.913  # K100
$ws.Cells.Item(100, 13).Value = -3563.913  # M100
$ws.Cells.Item(113, 8).Value = 2236.889  # H113
$ws.Cells.Item(113, 9).Value = 2110.6667  # I113
$ws.Cells.Item(113, 11).Value = 2110.6667  # K113
$ws.Cells.Item(113, 13).Value = 59.33329999999978  # M113
$ws.Cells.Item(122, 8).Value = 3705.5715  # H122
$ws.Cells.Item(122, 9).Value = 2435  # I122
$ws.Cells.Item(122, 11).Value = 7305  # K122
$ws.Cells.Item(122, 13).Value = -4855  # M122
$ws.Cells.Item(126, 8).Value = 4251  # H126
$ws.Cells.Item(126, 9).Value = 3797.1428  # I126
$ws.Cells.Item(126, 11).Value = 11391.4284  # K126
$ws.Cells.Item(126, 13).Value = -8921.428400000001  # M126
$ws.Cells.Item(132, 8).Value = 3705864  # H132
$ws.Cells.Item(132, 9).Value = 5557317.5  # I132
$ws.Cells.Item(132, 10).Value = 2956.7  # J132
$ws.Cells.Item(132, 11).Value = 16671952.5  # K132
$ws.Cells.Item(132, 12).Value = 8870.099999999999  # L132
$ws.Cells.Item(132, 13).Value = -16669422.5  # M132
$ws.Cells.Item(132, 14).Value = -13930.1  # N132
$ws.Cells.Item(136, 8).Value = 3369551  # H136
$ws.Cells.Item(136, 9).Value = 3970384.2  # I136
$ws.Cells.Item(136, 10).Value = 4884.2  # J136
$ws.Cells.Item(136, 11).Value = 11911152.6  # K136
$ws.Cells.Item(136, 12).Value = 14652.6  # L136
$ws.Cells.Item(136, 13).Value = -11908602.6  # M136
$ws.Cells.Item(136, 14).Value = -19752.6  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 16027  # H52
$ws.Cells.Item(52, 9).Value = 17520.5  # I52
$ws.Cells.Item(52, 10).Value = 15031.333  # J52
$ws.Cells.Item(52, 11).Value = 17520.5  # K52
$ws.Cells.Item(52, 12).Value = 15031.333  # L52
$ws.Cells.Item(52, 13).Value = -17294.5  # M52
$ws.Cells.Item(52, 14).Value = -15483.333  # N52
$ws.Cells.Item(64, 8).Value = 63999  # H64
$ws.Cells.Item(64, 10).Value = 69999  # J64
$ws.Cells.Item(64, 12).Value = 69999  # L64
$ws.Cells.Item(64, 14).Value = -70495  # N64
$ws.Cells.Item(67, 8).Value = 63999  # H67
$ws.Cells.Item(67, 10).Value = 69999  # J67
$ws.Cells.Item(67, 12).Value = 69999  # L67
$ws.Cells.Item(67, 14).Value = -71715  # N67
$ws.Cells.Item(81, 8).Value = 2642.2  # H81
$ws.Cells.Item(81, 9).Value = 1811  # I81
$ws.Cells.Item(81, 11).Value = 3622  # K81
$ws.Cells.Item(81, 13).Value = -2561  # M81
$ws.Cells.Item(84, 8).Value = 2642.2  # H84
$ws.Cells.Item(84, 9).Value = 1811  # I84
$ws.Cells.Item(84, 11).Value = 18110  # K84
$ws.Cells.Item(84, 13).Value = -12806  # M84
$ws.Cells.Item(100, 8).Value = 1165.2727  # H100
$ws.Cells.Item(100, 9).Value = 1081.9  # I100
$ws.Cells.Item(100, 11).Value = 2163.8  # K100
$ws.Cells.Item(100, 13).Value = -1622.8  # M100
$ws.Cells.Item(107, 8).Value = 1846.7391  # H107
$ws.Cells.Item(107, 9).Value = 1922.2727  # I107
$ws.Cells.Item(107, 10).Value = 1777.5  # J107
$ws.Cells.Item(107, 11).Value = 5766.8181  # K107
$ws.Cells.Item(107, 12).Value = 5332.5  # L107
$ws.Cells.Item(107, 13).Value = -3846.8181  # M107
$ws.Cells.Item(107, 14).Value = -9172.5  # N107
$ws.Cells.Item(113, 8).Value = 486.08334  # H113
$ws.Cells.Item(113, 9).Value = 191.75  # I113
$ws.Cells.Item(113, 11).Value = 575.25  # K113
$ws.Cells.Item(113, 13).Value = 1594.75  # M113
$ws.Cells.Item(122, 8).Value = 6478.4  # H122
$ws.Cells.Item(122, 9).Value = 3885.625  # I122
$ws.Cells.Item(122, 10).Value = 16849.5  # J122
$ws.Cells.Item(122, 11).Value = 11656.875  # K122
$ws.Cells.Item(122, 12).Value = 50548.5  # L122
$ws.Cells.Item(122, 13).Value = -9206.875  # M122
$ws.Cells.Item(122, 14).Value = -55448.5  # N122
$ws.Cells.Item(126, 8).Value = 4234.4707  # H126
$ws.Cells.Item(126, 9).Value = 5371.4546  # I126
$ws.Cells.Item(126, 10).Value = 2150  # J126
$ws.Cells.Item(126, 11).Value = 16114.3638  # K126
$ws.Cells.Item(126, 12).Value = 6450  # L126
$ws.Cells.Item(126, 13).Value = -13644.3638  # M126
$ws.Cells.Item(126, 14).Value = -11390  # N126
$ws.Cells.Item(136, 8).Value = 928.9259  # H136
$ws.Cells.Item(136, 9).Value = 928.9259  # I136
$ws.Cells.Item(136, 11).Value = 2786.7777  # K136
$ws.Cells.Item(136, 13).Value = -236.7776999999996  # M136
